$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "26.857.34"
$ws.Range("E2").Value = "  +0.02%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.638.50"
$ws.Range("E3").Value = "  -0.23%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.62%  "

# Row 5: BNB
$ws.Range("D5").Value = "'217.03"
$ws.Range("E5").Value = "  -0.75%  "

# Row 6: XRP
$ws.Range("D6").Value = "'0.508"
$ws.Range("E6").Value = "  +2.15%  "

# Row 7: USDC
$ws.Range("E7").Value = "  -0.62%  "

# Row 8: Cardano
$ws.Range("D8").Value = "'0.253"
$ws.Range("E8").Value = "  +1.25%  "

# Row 9: Dogecoin
$ws.Range("D9").Value = "'0.0623"
$ws.Range("E9").Value = "  +0.20%  "

# Row 10: Solana
$ws.Range("E10").Value = "  +3.23%  "

# Row 11: TRON
$ws.Range("D11").Value = "'0.0846"
$ws.Range("E11").Value = "  +0.12%  "

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.868.23"
$ws.Range("E12").Value = "  -0.20%  "

# Row 13: WrappedEther
$ws.Range("D13").Value = "1.641.16"
$ws.Range("E13").Value = "  -0.01%  "

# Row 14: Polkadot
$ws.Range("E14").Value = "  -1.03%  "

# Row 15: Polygon
$ws.Range("E15").Value = "  +0.54%  "

# Row 16: Litecoin
$ws.Range("D16").Value = "'67.11"
$ws.Range("E16").Value = "  +2.64%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "26.854.89"
$ws.Range("E17").Value = "  -0.03%  "

# Row 18: ShibaInu
$ws.Range("D18").Value = "0.0₃0730"
$ws.Range("E18").Value = "  +0.16%  "

# Row 19: BitcoinCash
$ws.Range("D19").Value = "'218.28"
$ws.Range("E19").Value = "  +1.15%  "

# Row 20: Dai
$ws.Range("E20").Value = "  -0.56%  "

# Row 21: Chainlink
$ws.Range("E21").Value = "  +1.39%  "

# Row 22: Uniswap
$ws.Range("E22").Value = "  +0.79%  "

# Row 23: Toncoin
$ws.Range("E23").Value = "  +2.88%  "

# Row 24: Avalanche
$ws.Range("D24").Value = "'9.16"
$ws.Range("E24").Value = "  -0.33%  "

# Row 25: Monero
$ws.Range("D25").Value = "'147.25"
$ws.Range("E25").Value = "  -0.28%  "

# Row 26: BinanceUSD
$ws.Range("E26").Value = "  -0.68%  "

# Row 27: Stellar
$ws.Range("E27").Value = "  +0.37%  "

# Row 28: Cosmos
$ws.Range("E28").Value = "  +0.41%  "

# Row 29: EthereumClassic
$ws.Range("D29").Value = "'15.75"
$ws.Range("E29").Value = "  +0.26%  "

# Row 30: Hedera
$ws.Range("E30").Value = "  -1.16%  "

# Row 31: PancakeSwap
$ws.Range("D31").Value = "'1.18"
$ws.Range("E31").Value = "  -1.38%  "

# Row 32: Filecoin
$ws.Range("D32").Value = "'3.33"
$ws.Range("E32").Value = "  -1.20%  "

# Row 33: InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +0.14%  "

# Row 34: LidoDAOToken
$ws.Range("E34").Value = "  +1.31%  "

# Row 35: Maker
$ws.Range("D35").Value = "1.265.78"
$ws.Range("E35").Value = "  -1.35%  "

# Row 36: HuobiToken
$ws.Range("E36").Value = "  -0.35%  "

# Row 37: VeChain
$ws.Range("E37").Value = "  +2.19%  "

# Row 38 & 39: ImmutableX and ARBITRUM swap places
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.533"
$ws.Range("E38").Value = "  +0.04%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'0.835"
$ws.Range("E39").Value = "  +2.04%  "

# Row 40: PaxDollar
$ws.Range("E40").Value = "  -0.60%  "

# Row 41: TrustWalletToken
$ws.Range("E41").Value = "  +0.34%  "

# Row 42: FraxShare
$ws.Range("E42").Value = "  +0.54%  "

# Row 43: RocketPoolETH
$ws.Range("D43").Value = "1.779.53"
$ws.Range("E43").Value = "  -0.19%  "

# Row 44: Aave
$ws.Range("D44").Value = "'61.96"
$ws.Range("E44").Value = "  +1.67%  "

# Row 45: MXToken
$ws.Range("E45").Value = "  -0.12%  "

# Row 46: Quant
$ws.Range("D46").Value = "'91.83"
$ws.Range("E46").Value = "  -0.90%  "

# Row 47: RenderToken
$ws.Range("E47").Value = "  -0.52%  "

# Row 48: BabyDogeCoin
$ws.Range("E48").Value = "  +0.07%  "

# Row 49: Cronos
$ws.Range("E49").Value = "  -0.78%  "

# Row 50: EnergySwap
$ws.Range("D50").Value = "'7.67"
$ws.Range("E50").Value = "  +1.12%  "

# Row 51: Algorand
$ws.Range("E51").Value = "  -0.67%  "
